$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Huy's row (row 5): change the data structure / technique notes,
# and clear the leftover placeholder text in F5, then add the file report
# text to E5.
$ws.Range("E5").Value = "HashMap"
$ws.Range("C5").Value = "Hashtable, ADT, Analysis of Algorithms"
$ws.Range("F5").Value = ""

# Reflect the final cell selection in the saved view.
$ws.Range("C5").Select()
